$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.300.78"
$ws.Range("E2").Value = "  +3.65%  "

$ws.Range("D3").Value = "2.498.59"
$ws.Range("E3").Value = "  +2.71%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.80%  "

$ws.Range("E7").Value = "  +1.63%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  +2.56%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.90"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.44%  "

$ws.Range("E11").Value = "  +1.68%  "

$ws.Range("E12").Value = "  +1.17%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.57"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.62%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.22"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.12%  "

$ws.Range("D15").Value = "2.888.61"
$ws.Range("E15").Value = "  +2.95%  "

$ws.Range("D16").Value = "2.497.28"
$ws.Range("E16").Value = "  +2.91%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.858"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.17%  "

$ws.Range("D18").Value = "47.282.23"
$ws.Range("E18").Value = "  +4.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.55%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.90%  "

$ws.Range("D21").Value = "0.0₃0949"
$ws.Range("E21").Value = "  +2.24%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.53%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.47%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "250.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.71%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.56%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.85%  "

$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("E28").Value = "  +4.93%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.79"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.06%  "

$ws.Range("E31").Value = "  +6.42%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.55"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.33%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.07"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.30%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.46"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.77%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0794"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.56%  "

$ws.Range("E36").Value = "  +0.28%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.78"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.92%  "

$ws.Range("E38").Value = "  +5.48%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.46%  "

$ws.Range("E40").Value = "  +1.64%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "122.56"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.30%  "

$ws.Range("E42").Value = "  -1.63%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.55"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.49%  "

$ws.Range("E44").Value = "  +2.75%  "

$ws.Range("D45").Value = "1.985.27"
$ws.Range("E45").Value = "  +1.49%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.05"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.09"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.29%  "

$ws.Range("E48").Value = "  -1.81%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.08"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.48%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +11.12%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.20%  "
